$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.815.29"
$ws.Range("E2").Value = "  -2.97%  "

# Row 3
$ws.Range("D3").Value = "1.614.44"
$ws.Range("E3").Value = "  -3.37%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").Value = "'306.66"
$ws.Range("E6").Value = "  -1.98%  "

# Row 7
$ws.Range("D7").Value = "'0.3893"
$ws.Range("E7").Value = "  -0.37%  "

# Row 8
$ws.Range("D8").Value = "'0.3796"
$ws.Range("E8").Value = "  -3.67%  "

# Row 9
$ws.Range("D9").Value = "'1.001"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").Value = "'1.350"
$ws.Range("E10").Value = "  -3.63%  "

# Row 11
$ws.Range("D11").Value = "'48.53"
$ws.Range("E11").Value = "  -5.48%  "

# Row 12
$ws.Range("E12").Value = "  -2.50%  "

# Row 13
$ws.Range("D13").Value = "'23.75"
$ws.Range("E13").Value = "  -6.10%  "

# Row 14
$ws.Range("D14").Value = "'6.989"
$ws.Range("E14").Value = "  -4.57%  "

# Row 15
$ws.Range("D15").Value = "'0.00001270"
$ws.Range("E15").Value = "  -3.98%  "

# Row 16
$ws.Range("D16").Value = "'7.422"
$ws.Range("E16").Value = "  -3.88%  "

# Row 17
$ws.Range("D17").Value = "1.612.03"
$ws.Range("E17").Value = "  -3.56%  "

# Row 18
$ws.Range("D18").Value = "'93.05"
$ws.Range("E18").Value = "  -0.30%  "

# Row 19
$ws.Range("D19").Value = "'0.06906"
$ws.Range("E19").Value = "  -1.51%  "

# Row 20
$ws.Range("D20").Value = "'19.92"
$ws.Range("E20").Value = "  -5.01%  "

# Row 21
$ws.Range("E21").Value = "  -3.93%  "

# Row 22
$ws.Range("E22").Value = "  -0.11%  "

# Row 23
$ws.Range("D23").Value = "'13.37"
$ws.Range("E23").Value = "  -4.06%  "

# Row 24
$ws.Range("D24").Value = "23.825.32"
$ws.Range("E24").Value = "  -2.94%  "

# Row 25
$ws.Range("D25").Value = "'2.431"
$ws.Range("E25").Value = "  +3.14%  "

# Row 26
$ws.Range("D26").Value = "'2.795"
$ws.Range("E26").Value = "  +1.94%  "

# Row 27
$ws.Range("D27").Value = "'22.06"
$ws.Range("E27").Value = "  -4.72%  "

# Row 28
$ws.Range("D28").Value = "'157.00"
$ws.Range("E28").Value = "  -2.01%  "

# Row 29
$ws.Range("D29").Value = "'138.98"
$ws.Range("E29").Value = "  -5.29%  "

# Row 30
$ws.Range("E30").Value = "  -10.57%  "

# Row 31
$ws.Range("D31").Value = "'7.720"
$ws.Range("E31").Value = "  -8.07%  "

# Row 32
$ws.Range("D32").Value = "'2.484"
$ws.Range("E32").Value = "  -0.94%  "

# Row 33
$ws.Range("D33").Value = "1.790.78"
$ws.Range("E33").Value = "  -3.60%  "

# Row 34
$ws.Range("D34").Value = "'0.08072"
$ws.Range("E34").Value = "  -3.23%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.9582"
$ws.Range("E35").Value = "  -3.23%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02862"
$ws.Range("E36").Value = "  -5.69%  "

# Row 37
$ws.Range("D37").Value = "'6.526"
$ws.Range("E37").Value = "  -6.84%  "

# Row 38
$ws.Range("D38").Value = "'0.2640"
$ws.Range("E38").Value = "  -6.04%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.09121"
$ws.Range("E39").Value = "  -3.57%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'10.37"
$ws.Range("E40").Value = "  +0.51%  "

# Row 41
$ws.Range("D41").Value = "'13.34"
$ws.Range("E41").Value = "  -1.77%  "

# Row 42
$ws.Range("D42").Value = "'1.422"
$ws.Range("E42").Value = "  -5.97%  "

# Row 43
$ws.Range("D43").Value = "'0.7415"
$ws.Range("E43").Value = "  -6.20%  "

# Row 44
$ws.Range("D44").Value = "'15.73"
$ws.Range("E44").Value = "  -4.21%  "

# Row 45
$ws.Range("D45").Value = "'0.6799"
$ws.Range("E45").Value = "  -4.51%  "

# Row 46
$ws.Range("D46").Value = "'2.430"
$ws.Range("E46").Value = "  -4.56%  "

# Row 47
$ws.Range("D47").Value = "'4.053"
$ws.Range("E47").Value = "  -2.94%  "

# Row 48
$ws.Range("E48").Value = "  +0.07%  "

# Row 49
$ws.Range("D49").Value = "'0.08205"
$ws.Range("E49").Value = "  -4.80%  "

# Row 50
$ws.Range("D50").Value = "'132.25"

# Row 51
$ws.Range("D51").Value = "'1.183"
$ws.Range("E51").Value = "  -10.57%  "
